$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.848.14"
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = "'1.564.74"
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = "'206.38"
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = "'1.786.13"
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = "'1.564.62"
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = "'3.75"
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("E15").Value = '  -0.96%  '
$ws.Range("D16").Value = "'26.911.42"
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("D17").Value = "'61.73"
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = "'215.54"
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = "'7.34"
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = "'9.36"
$ws.Range("E23").Value = '  -2.77%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'152.09"
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("D26").Value = "'6.70"
$ws.Range("E26").Value = '  -1.80%  '
$ws.Range("D27").Value = "'14.90"
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").Value = "'3.15"
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("D33").Value = "'1.393.26"
$ws.Range("E33").Value = '  +2.32%  '
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("D36").Value = "'2.29"
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("D37").Value = "'0.942"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  +2.98%  '
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").Value = "'1.79"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("E45").Value = '  +1.74%  '
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = "'1.699.16"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = "'85.37"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = "'0.0₇0971"
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.0494"
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = "'0.0947"
$ws.Range("E51").Value = '  -0.86%  '
